$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4054.4443
$ws.Range("I74").Value = 3496.6667
$ws.Range("K74").Value = 3496.6667
$ws.Range("M74").Value = -2560.6667
$ws.Range("H77").Value = 4054.4443
$ws.Range("I77").Value = 3496.6667
$ws.Range("K77").Value = 17483.3335
$ws.Range("M77").Value = -12803.3335
$ws.Range("H93").Value = 19000
$ws.Range("I93").Value = 19000
$ws.Range("K93").Value = 19000
$ws.Range("M93").Value = -16504
$ws.Range("H100").Value = 13334858
$ws.Range("I100").Value = 15152520
$ws.Range("J100").Value = 5337.3335
$ws.Range("K100").Value = 15152520
$ws.Range("L100").Value = 5337.3335
$ws.Range("M100").Value = -15151979
$ws.Range("N100").Value = -6419.3335
$ws.Range("H107").Value = 726.3333
$ws.Range("I107").Value = 343.05884
$ws.Range("J107").Value = 1657.1428
$ws.Range("K107").Value = 343.05884
$ws.Range("L107").Value = 1657.1428
$ws.Range("M107").Value = 1576.94116
$ws.Range("N107").Value = -5497.1428
$ws.Range("H133").Value = 61800
$ws.Range("J133").Value = 61800
$ws.Range("L133").Value = 61800
$ws.Range("N133").Value = -71920
$ws.Range("H137").Value = 1284.68
$ws.Range("I137").Value = 1389.2354
$ws.Range("J137").Value = 1062.5
$ws.Range("K137").Value = 4167.706200000001
$ws.Range("L137").Value = 3187.5
$ws.Range("M137").Value = -1617.706200000001
$ws.Range("N137").Value = -8287.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4325
$ws.Range("I63").Value = 3100
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 3100
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = -2414
$ws.Range("N63").Value = -9372
$ws.Range("H66").Value = 4325
$ws.Range("I66").Value = 3100
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 15500
$ws.Range("L66").Value = 40000
$ws.Range("M66").Value = -12068
$ws.Range("N66").Value = -46864
$ws.Range("H74").Value = 784.2632
$ws.Range("I74").Value = 644.3333
$ws.Range("J74").Value = 1127.7273
$ws.Range("K74").Value = 644.3333
$ws.Range("L74").Value = 1127.7273
$ws.Range("M74").Value = 229.6667
$ws.Range("N74").Value = -2875.7273
$ws.Range("H77").Value = 784.2632
$ws.Range("I77").Value = 644.3333
$ws.Range("J77").Value = 1127.7273
$ws.Range("K77").Value = 3221.6665
$ws.Range("L77").Value = 5638.636500000001
$ws.Range("M77").Value = 1146.3335
$ws.Range("N77").Value = -14374.6365
$ws.Range("H97").Value = 753.73334
$ws.Range("I97").Value = 593.2857
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 593.2857
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -97.28570000000002
$ws.Range("N97").Value = -3992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4091.4167
$ws.Range("I105").Value = 3670
$ws.Range("J105").Value = 4512.8335
$ws.Range("K105").Value = 3670
$ws.Range("L105").Value = 4512.8335
$ws.Range("M105").Value = -1923
$ws.Range("N105").Value = -8006.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1047.6
$ws.Range("I58").Value = 955.7143
$ws.Range("J58").Value = 1164.5454
$ws.Range("K58").Value = 955.7143
$ws.Range("L58").Value = 1164.5454
$ws.Range("M58").Value = -752.7143
$ws.Range("N58").Value = -1570.5454
$ws.Range("H136").Value = 1047.6
$ws.Range("I136").Value = 955.7143
$ws.Range("J136").Value = 1164.5454
$ws.Range("K136").Value = 2867.1429
$ws.Range("L136").Value = 3493.6362
$ws.Range("M136").Value = -317.1428999999998
$ws.Range("N136").Value = -8593.636200000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 7526.5
$ws.Range("J74").Value = 8571.286
$ws.Range("L74").Value = 25713.858
$ws.Range("N74").Value = -27835.858
$ws.Range("H77").Value = 7526.5
$ws.Range("J77").Value = 8571.286
$ws.Range("L77").Value = 77141.57399999999
$ws.Range("N77").Value = -87749.57399999999
$ws.Range("H102").Value = 5487.6665
$ws.Range("H120").Value = 12260
$ws.Range("J120").Value = 19500
$ws.Range("L120").Value = 58500
$ws.Range("N120").Value = -68176
$ws.Range("H131").Value = 6411211
$ws.Range("I131").Value = 303.3
$ws.Range("J131").Value = 7353991.5
$ws.Range("K131").Value = 909.9000000000001
$ws.Range("L131").Value = 22061974.5
$ws.Range("M131").Value = 4130.1
$ws.Range("N131").Value = -22072054.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2272.2727
$ws.Range("I113").Value = 1882.5
$ws.Range("J113").Value = 2740
$ws.Range("K113").Value = 1882.5
$ws.Range("L113").Value = 2740
$ws.Range("M113").Value = 287.5
$ws.Range("N113").Value = -7080
$ws.Range("H126").Value = 11120933
$ws.Range("I126").Value = 10450.75
$ws.Range("J126").Value = 55562864
$ws.Range("K126").Value = 31352.25
$ws.Range("L126").Value = 166688592
$ws.Range("M126").Value = -28882.25
$ws.Range("N126").Value = -166693532
$ws.Range("H133").Value = 58400
$ws.Range("J133").Value = 58400
$ws.Range("L133").Value = 58400
$ws.Range("N133").Value = -68520

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2491.3044
$ws.Range("J100").Value = 2838.889
$ws.Range("L100").Value = 2838.889
$ws.Range("N100").Value = -3920.889
$ws.Range("H133").Value = 111952.6
$ws.Range("J133").Value = 111952.6
$ws.Range("L133").Value = 111952.6
$ws.Range("N133").Value = -117012.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 27000
$ws.Range("J110").Value = 27000
$ws.Range("L110").Value = 27000
$ws.Range("N110").Value = -35180
$ws.Range("H122").Value = 1701.8
$ws.Range("I122").Value = 2004
$ws.Range("J122").Value = 1626.25
$ws.Range("K122").Value = 6012
$ws.Range("L122").Value = 4878.75
$ws.Range("M122").Value = -3562
$ws.Range("N122").Value = -9778.75
$ws.Range("H132").Value = 56252224
$ws.Range("I132").Value = 93751240
$ws.Range("J132").Value = 3697.0625
$ws.Range("K132").Value = 281253720
$ws.Range("L132").Value = 11091.1875
$ws.Range("M132").Value = -281251190
$ws.Range("N132").Value = -16151.1875
$ws.Range("H133").Value = 42107.5
$ws.Range("J133").Value = 42107.5
$ws.Range("L133").Value = 42107.5
$ws.Range("N133").Value = -52227.5
$ws.Range("H136").Value = 52994.21
$ws.Range("I136").Value = 62849.375
$ws.Range("J136").Value = 433.33334
$ws.Range("K136").Value = 188548.125
$ws.Range("L136").Value = 1300.00002
$ws.Range("M136").Value = -185998.125
$ws.Range("N136").Value = -6400.000019999999
